# Daily attendance processing - 2025-12-24 18:40:42
# Reorders the "Recorded By" names in column G so the first contributor in
# the comma-separated list is moved to the end (e.g.
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com",
# "system, backup@backdoor.com, System" -> "backup@backdoor.com, System, system").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2, 3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 28, 29, 32, 36, 38, 39, 40, 41, 44, 45, 46, 47, 48, 50, 52, 54, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 99, 101, 109, 110, 111, 112, 116, 118, 125, 127, 135, 136, 137, 138, 142, 144, 151, 153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $current = [string]$cell.Value2
    $parts = $current -split ", "
    $lastIndex = $parts.Length - 1
    # Left-rotate: move the first name to the end of the list.
    $reordered = @($parts[1..$lastIndex]) + @($parts[0])
    $cell.Value = [string]::Join(", ", $reordered)
}
